# Apply a green highlight (RGB 00FF00) to the text run in three
# "dataset name" table cells across the presentation, matching the
# commit "update list of datasets done in the pptx file".

$p = $ppt.ActivePresentation

# Map of slide index -> row/col of the target table cell whose text
# is the dataset name that should be marked done (highlighted green).
$targets = @(
    @{ Slide = 2; Row = 7; Col = 1; Text = "Civil Liberties index" },
    @{ Slide = 3; Row = 7; Col = 1; Text = "Economic Group" },
    @{ Slide = 4; Row = 5; Col = 1; Text = "Geographical regions" }
)

foreach ($t in $targets) {
    $slide = $p.Slides.Item($t.Slide)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $cell = $shape.Table.Cell($t.Row, $t.Col)
            $textRange = $cell.Shape.TextFrame.TextRange
            if ($textRange.Text -eq $t.Text) {
                $textRange.Font.Highlight.RGB = 65280
            }
        }
    }
}
